$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-12 (symptom_group label + 5 numeric columns B-F)
$data = @(
    @{ Row = 2;  Label = "Diseases (patient-stated)"; Values = @(6.46, 3.18, 4.89, 6.07, 4.87) }
    @{ Row = 3;  Label = "Injuries & adverse effects"; Values = @(10.58, 4.34, 6.67, 8.99, 6.73) }
    @{ Row = 4;  Label = "Other"; Values = @(6.88, 5.4, 7.48, 6.42, 5.78) }
    @{ Row = 5;  Label = "Symptom – Circulatory"; Values = @(10.54, 6.91, 8.59, 10.09, 8.62) }
    @{ Row = 6;  Label = "Symptom – Digestive"; Values = @(11.83, 6.71, 11.7, 11.83, 10.94) }
    @{ Row = 7;  Label = "Symptom – General"; Values = @(4.72, 4.44, 5.48, 5.01, 5.16) }
    @{ Row = 8;  Label = "Symptom – Genitourinary"; Values = @(5.83, 4.14, 5.04, 5.7, 5.3) }
    @{ Row = 9;  Label = "Symptom – Nervous"; Values = @(10.45, 11.7, 11.26, 11.47, 12.68) }
    @{ Row = 10; Label = "Symptom – Respiratory"; Values = @(27.39, 50.28, 34.22, 29.94, 36.58) }
    @{ Row = 11; Label = "Symptom – Skin/Hair/Nails"; Values = @(2.65, 1.87, 2.67, 2.41, 2.02) }
    @{ Row = 12; Label = "Uncodable/Unknown"; Values = @(2.67, 1.06, 2, 2.07, 1.33) }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.Label
    for ($i = 0; $i -lt $entry.Values.Length; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $entry.Values[$i]
    }
}
